$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 88 (pushes existing rows 88..122 down to 89..123)
$ws.Rows.Item(88).Insert()

# Populate the new row 88 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,O,R carry the same constant values as the rest of the
# "Comercializadora del Agro de Limarí" / "Poroto granado" block (copied from row 89,
# which now holds what used to be row 88's data).
$ws.Cells.Item(88, 1).Value = 2
$ws.Cells.Item(88, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(88, 3).Value = "Coquimbo"
$ws.Cells.Item(88, 4).Value = 45029
$ws.Cells.Item(88, 5).Value = 4
$ws.Cells.Item(88, 6).Value = 100112030
$ws.Cells.Item(88, 7).Value = "Poroto granado"
$ws.Cells.Item(88, 8).Value = "Sin especificar"
$ws.Cells.Item(88, 9).Value = "Primera"
$ws.Cells.Item(88, 10).Value = 500
$ws.Cells.Item(88, 11).Value = 30000
$ws.Cells.Item(88, 12).Value = 31000
$ws.Cells.Item(88, 13).Value = 30500
$ws.Cells.Item(88, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(88, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(88, 16).Value = 1220
$ws.Cells.Item(88, 17).Value = 25
$ws.Cells.Item(88, 18).Value = "Hortaliza"
